$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.1371
$ws.Range("C6").Value = -11.7912
$ws.Range("A14").Value = -21.86449999999999
$ws.Range("B15").Value = 4.828099999999999
$ws.Range("A16").Value = -22.20530000000001
$ws.Range("C18").Value = -12.1419
$ws.Range("C19").Value = -12.87050000000001
$ws.Range("A21").Value = -22.0425
$ws.Range("B21").Value = 5.244399999999997
$ws.Range("B22").Value = 9.967800000000004
$ws.Range("A23").Value = -20.19729999999998
$ws.Range("B24").Value = 5.520000000000003
$ws.Range("A25").Value = -21.67829999999998
$ws.Range("A26").Value = -21.08889999999997
$ws.Range("B27").Value = 5.784100000000006
$ws.Range("B28").Value = 5.838199999999998
$ws.Range("A29").Value = -20.92089999999996
$ws.Range("C35").Value = -13.65780000000002
$ws.Range("B36").Value = 8.523600000000005
$ws.Range("B39").Value = 8.514900000000006
$ws.Range("A40").Value = -20.2717
$ws.Range("C44").Value = -13.0732
$ws.Range("B45").Value = 4.460900000000004
$ws.Range("C47").Value = -12.0772
$ws.Range("B48").Value = 5.090600000000006
$ws.Range("B49").Value = 5.418999999999993
$ws.Range("C50").Value = -13.62669999999998
$ws.Range("C51").Value = -11.8974
$ws.Range("B52").Value = 5.358800000000005
$ws.Range("C52").Value = -11.2093
$ws.Range("A53").Value = -20.19619999999998
$ws.Range("B53").Value = 9.6425
$ws.Range("B54").Value = 4.828200000000002
$ws.Range("C55").Value = -13.5501
$ws.Range("A57").Value = -20.20779999999999
$ws.Range("B57").Value = 8.377900000000004
$ws.Range("C57").Value = -12.36960000000001
$ws.Range("C58").Value = -12.79969999999999
$ws.Range("A59").Value = -22.5829
$ws.Range("C64").Value = -10.44319999999999
$ws.Range("A65").Value = -21.85429999999999
$ws.Range("C66").Value = -11.15300000000001
$ws.Range("A69").Value = -21.61029999999999
$ws.Range("B70").Value = 5.0501
$ws.Range("B71").Value = 4.562799999999997
$ws.Range("A79").Value = -19.9508
$ws.Range("C80").Value = -13.16130000000001
$ws.Range("A83").Value = -21.82269999999999
$ws.Range("C83").Value = -12.2131
$ws.Range("B86").Value = 5.027600000000001
$ws.Range("B87").Value = 5.352999999999997
$ws.Range("B89").Value = 4.476699999999997
$ws.Range("A91").Value = -20.11889999999998
$ws.Range("C92").Value = -10.3361
$ws.Range("A93").Value = -21.24170000000002
$ws.Range("C94").Value = -10.7918
$ws.Range("C96").Value = -10.1739
$ws.Range("C97").Value = -11.1003
$ws.Range("A100").Value = -22.1281
$ws.Range("B101").Value = 5.776000000000002
$ws.Range("C101").Value = -12.2974
$ws.Range("A103").Value = -21.932
